# ---------------------------------------------------------------------------
# Applies the "add sth about sports meeting" commit to 入学志願理由書.docx
#
# Strategy: the COM-interop runtime backing this script merges every run
# that a text replacement touches together with any sibling runs that follow
# it in the same paragraph. To avoid corrupting not-yet-edited text we apply
# every text replacement from the END of the document towards the START.
# After the text content is correct, we make a second pass that re-applies
# precise character-level formatting (Bold on/off, and a couple of runs
# whose font was set "fresh" rather than inherited) using narrow Range
# objects, which this runtime is able to splice in without disturbing
# neighbouring runs.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ===========================================================================
# STEP 1 (bottom of doc first): paragraph 3 - merge "...いたしま" + "した。"
# into "...いたしました。" - this also removes the old _GoBack bookmark that
# used to sit at that exact seam.
# ===========================================================================
$null = $d.Content.Find.Execute(
    "貴校を志願いたしまし", $true, $false, $false, $false, $false,
    $true, 1, $false, "貴校を志願いたしまし", 1)

$null = $d.Content.Find.Execute(
    "いたしまし", $true, $false, $false, $false, $false,
    $true, 1, $false, "いたしました", 1)

# ===========================================================================
# STEP 2: paragraph 2 - second half of the touched region:
#   "運動会を応援の時、生徒たちの旺盛な活力と"
#   -> "運動会を応援の時、「玉入れ」の試合は何度も継続的な競争、子供たちは
#       まだ途中をあきらめなかった、すべて最終的な終わりに従う。各試合の
#       結果は、子供たち自身が発表するものです。生徒たちの旺盛な活力と"
# ===========================================================================
$null = $d.Content.Find.Execute(
    "運動会を応援の時、生徒たちの旺盛な活力と", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "運動会を応援の時、「玉入れ」の試合は何度も継続的な競争、子供たちはまだ途中をあきらめなかった、すべて最終的な終わりに従う。各試合の結果は、子供たち自身が発表するものです。生徒たちの旺盛な活力と",
    1)

# ===========================================================================
# STEP 3: paragraph 2 - first half of the touched region:
#   "になる信じております。貴校の学校公開で訪問した際、児童挨拶など※※※※※※※※※※※※※※"
#   -> "になる信じております。貴校の学校公開で授業を見学した際、先生方のとても
#       パワフルなご指導に関心しました。また、内容もレベルが高く、児童達も
#       熱心に耳を傾けていたことが印象的でした"
# ===========================================================================
$null = $d.Content.Find.Execute(
    "になる信じております。貴校の学校公開で訪問した際、児童挨拶など※※※※※※※※※※※※※※",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "になる信じております。貴校の学校公開で授業を見学した際、先生方のとてもパワフルなご指導に関心しました。また、内容もレベルが高く、児童達も熱心に耳を傾けていたことが印象的でした",
    1)

Write-Output "text pass done"
